$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.503.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.24%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.524.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.06%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '615.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.14%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.19%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.524.03'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 10
$ws.Range('E10').Value = '  -0.77%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.14%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.426'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.86%  '

# Row 13
$ws.Range('E13').Value = '  -0.47%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.63%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.118.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.521.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.466.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.26%  '

# Row 18
$ws.Range('E18').Value = '  +0.09%  '

# Row 19
$ws.Range('E19').Value = '  +0.22%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.65%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '444.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.67%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.91%  '

# Row 23
$ws.Range('E23').Value = '  -2.67%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.58%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000131'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.71%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.662.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.22%  '

# Row 27
$ws.Range('E27').Value = '  +0.12%  '

# Row 28
$ws.Range('E28').Value = '  -1.76%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.56%  '

# Row 30
$ws.Range('E30').Value = '  -0.95%  '

# Row 31
$ws.Range('E31').Value = '  -4.94%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.02%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.164'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.31%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.22%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.69%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.515.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '

# Row 37
$ws.Range('E37').Value = '  -2.81%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.39%  '

# Row 39
$ws.Range('E39').Value = '  +0.05%  '

# Row 40
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '177.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.62%  '

# Row 41
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.02%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0887'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.03%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.06%  '

# Row 44
$ws.Range('E44').Value = '  -3.38%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.883'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.25%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.99%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.35%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.30%  '

# Row 49
$ws.Range('E49').Value = '  +2.74%  '

# Row 50
$ws.Range('E50').Value = '  -0.52%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.996'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.89%  '
